# Apply updated vm_pu values (380 kV case) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.06857615587679
$ws.Range("D2").Value = 1.07770985698997
$ws.Range("E2").Value = 1.063204140654213
$ws.Range("F2").Value = 1.083722961351984
$ws.Range("I2").Value = 1.045368216597133
$ws.Range("J2").Value = 1.073515261071374
$ws.Range("K2").Value = 1.080390442426415
$ws.Range("L2").Value = 1.065923413181704
$ws.Range("M2").Value = 1.08638783758129
$ws.Range("N2").Value = 1.075039775876017
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.070282543834819
$ws.Range("D3").Value = 1.079353904945833
$ws.Range("E3").Value = 1.064684570457173
$ws.Range("F3").Value = 1.085408002656506
$ws.Range("I3").Value = 1.045759729291559
$ws.Range("J3").Value = 1.074874907245315
$ws.Range("K3").Value = 1.081850220615813
$ws.Range("L3").Value = 1.067217138152889
$ws.Range("M3").Value = 1.087889663993227
$ws.Range("N3").Value = 1.076401352903478
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.071384231814859
$ws.Range("D4").Value = 1.080415557570573
$ws.Range("E4").Value = 1.065640453174503
$ws.Range("F4").Value = 1.08649622241468
$ws.Range("I4").Value = 1.04601060252483
$ws.Range("J4").Value = 1.075751880320371
$ws.Range("K4").Value = 1.082792156361305
$ws.Range("L4").Value = 1.068051700978114
$ws.Range("M4").Value = 1.088858854783705
$ws.Range("N4").Value = 1.077279571380891
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.07184680564474
$ws.Range("D5").Value = 1.080861371791513
$ws.Range("E5").Value = 1.066041824736313
$ws.Range("F5").Value = 1.086953215682903
$ws.Range("I5").Value = 1.04611548382185
$ws.Range("J5").Value = 1.076119897623422
$ws.Range("K5").Value = 1.083187525337855
$ws.Range("L5").Value = 1.068401946994975
$ws.Range("M5").Value = 1.089265693260853
$ws.Range("N5").Value = 1.077648111310725
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.071924440434828
$ws.Range("D6").Value = 1.080936196712166
$ws.Range("E6").Value = 1.066109188847647
$ws.Range("F6").Value = 1.087029918221363
$ws.Range("I6").Value = 1.04613305960247
$ws.Range("J6").Value = 1.076181650814722
$ws.Range("K6").Value = 1.083253873436721
$ws.Range("L6").Value = 1.068460719688663
$ws.Range("M6").Value = 1.089333967818873
$ws.Range("N6").Value = 1.07770995219864
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.071390414993717
$ws.Range("D7").Value = 1.080421516528397
$ws.Range("E7").Value = 1.065645818198842
$ws.Range("F7").Value = 1.086502330706973
$ws.Range("I7").Value = 1.046012006250796
$ws.Range("J7").Value = 1.075756800373302
$ws.Range("K7").Value = 1.082797441725359
$ws.Range("L7").Value = 1.068056383343391
$ws.Range("M7").Value = 1.088864293358375
$ws.Range("N7").Value = 1.077284498420862
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.06915335466437
$ws.Range("D8").Value = 1.078265924240641
$ws.Range("E8").Value = 1.063704890420958
$ws.Range("F8").Value = 1.084292874009708
$ws.Range("I8").Value = 1.045501041974179
$ws.Range("J8").Value = 1.073975347680123
$ws.Range("K8").Value = 1.080884334125436
$ws.Range("L8").Value = 1.066361170215214
$ws.Range("M8").Value = 1.08689593006817
$ws.Range("N8").Value = 1.075500515860524
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.065191952553679
$ws.Range("D9").Value = 1.074450453738555
$ws.Range("E9").Value = 1.06026853037189
$ws.Range("F9").Value = 1.080382795982733
$ws.Range("I9").Value = 1.044581645540955
$ws.Range("J9").Value = 1.070814224165494
$ws.Range("K9").Value = 1.077492500119106
$ws.Range("L9").Value = 1.06335392169161
$ws.Range("M9").Value = 1.083407078807398
$ws.Range("N9").Value = 1.072334903187917
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.062537146503681
$ws.Range("D10").Value = 1.071894609341644
$ws.Range("E10").Value = 1.057966082247328
$ws.Range("F10").Value = 1.077764072024835
$ws.Range("I10").Value = 1.043955708173213
$ws.Range("J10").Value = 1.068691381738865
$ws.Range("K10").Value = 1.075216679173049
$ws.Range("L10").Value = 1.061334991708064
$ws.Range("M10").Value = 1.0810668013633
$ws.Range("N10").Value = 1.070209046081549
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.061384113324555
$ws.Range("D11").Value = 1.070784841864867
$ws.Range("E11").Value = 1.056966213446159
$ws.Range("F11").Value = 1.076627119169839
$ws.Range("I11").Value = 1.043681535097239
$ws.Range("J11").Value = 1.06776836266352
$ws.Range("K11").Value = 1.074227607706772
$ws.Range("L11").Value = 1.06045729317896
$ws.Range("M11").Value = 1.080049869000018
$ws.Range("N11").Value = 1.069284716213246
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.060955286099508
$ws.Range("D12").Value = 1.070372149221565
$ws.Range("E12").Value = 1.056594370267983
$ws.Range("F12").Value = 1.076204334666489
$ws.Range("I12").Value = 1.043579219142487
$ws.Range("J12").Value = 1.067424926561356
$ws.Range("K12").Value = 1.073859664452387
$ws.Range("L12").Value = 1.060130740960965
$ws.Range("M12").Value = 1.079671584091073
$ws.Range("N12").Value = 1.068940792392424
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.061047295658307
$ws.Range("D13").Value = 1.070460694948078
$ws.Range("E13").Value = 1.056674152341822
$ws.Range("F13").Value = 1.076295044864471
$ws.Range("I13").Value = 1.043601187867376
$ws.Range("J13").Value = 1.067498621480943
$ws.Range("K13").Value = 1.07393861499868
$ws.Range("L13").Value = 1.060200811961739
$ws.Range("M13").Value = 1.079752752639363
$ws.Range("N13").Value = 1.069014591967249
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.061348677427408
$ws.Range("D14").Value = 1.070750738347732
$ws.Range("E14").Value = 1.056935485973316
$ws.Range("F14").Value = 1.076592181335714
$ws.Range("I14").Value = 1.043673087361717
$ws.Range("J14").Value = 1.067739986150444
$ws.Range("K14").Value = 1.074197204864621
$ws.Range("L14").Value = 1.060430311271721
$ws.Range("M14").Value = 1.080018611169722
$ws.Range("N14").Value = 1.069256299402266
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.061534296734732
$ws.Range("D15").Value = 1.070929380206435
$ws.Range("E15").Value = 1.057096442584704
$ws.Range("F15").Value = 1.076775194256146
$ws.Range("I15").Value = 1.043717323823472
$ws.Range("J15").Value = 1.067888620894544
$ws.Range("K15").Value = 1.074356456264786
$ws.Range("L15").Value = 1.060571642003476
$ws.Range("M15").Value = 1.080182341969765
$ws.Range("N15").Value = 1.069405145224751
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.062613594753084
$ws.Range("D16").Value = 1.071968195050867
$ws.Range("E16").Value = 1.058032378204946
$ws.Range("F16").Value = 1.077839462795623
$ws.Range("I16").Value = 1.043973837646765
$ws.Range("J16").Value = 1.068752558034196
$ws.Range("K16").Value = 1.075282243060681
$ws.Range("L16").Value = 1.061393167165855
$ws.Range("M16").Value = 1.08113421530056
$ws.Range("N16").Value = 1.070270309254236
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.063289665820089
$ws.Range("D17").Value = 1.07261898414504
$ws.Range("E17").Value = 1.058618682306001
$ws.Range("F17").Value = 1.078506229328382
$ws.Range("I17").Value = 1.04413389874561
$ws.Range("J17").Value = 1.069293452882404
$ws.Range("K17").Value = 1.075861984528186
$ws.Range("L17").Value = 1.061907545932369
$ws.Range("M17").Value = 1.081730333144322
$ws.Range("N17").Value = 1.07081197223517
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.063683671496067
$ws.Range("D18").Value = 1.072998283181778
$ws.Range("E18").Value = 1.058960385243426
$ws.Range("F18").Value = 1.078894851796125
$ws.Range("I18").Value = 1.044226957199257
$ws.Range("J18").Value = 1.069608580534457
$ws.Range("K18").Value = 1.076199788969892
$ws.Range("L18").Value = 1.062207238575595
$ws.Range("M18").Value = 1.082077694413752
$ws.Range("N18").Value = 1.071127547404636
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.06381796088069
$ws.Range("D19").Value = 1.073127564707013
$ws.Range("E19").Value = 1.059076850286779
$ws.Range("F19").Value = 1.079027313068101
$ws.Range("I19").Value = 1.04425863655423
$ws.Range("J19").Value = 1.069715969024229
$ws.Range("K19").Value = 1.076314912767742
$ws.Range("L19").Value = 1.062309369421399
$ws.Range("M19").Value = 1.082196077659078
$ws.Range("N19").Value = 1.07123508839838
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.063217164572275
$ws.Range("D20").Value = 1.072549191275149
$ws.Range("E20").Value = 1.058555806304641
$ws.Range("F20").Value = 1.07843472181257
$ws.Range("I20").Value = 1.044116757023002
$ws.Range("J20").Value = 1.06923545805292
$ws.Range("K20").Value = 1.075799819944816
$ws.Range("L20").Value = 1.061852392776156
$ws.Range("M20").Value = 1.081666411039881
$ws.Range("N20").Value = 1.070753895046376
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.061259943025219
$ws.Range("D21").Value = 1.070665341040359
$ws.Range("E21").Value = 1.056858542184804
$ws.Range("F21").Value = 1.076504695132292
$ws.Range("I21").Value = 1.043651927924582
$ws.Range("J21").Value = 1.06766892654226
$ws.Range("K21").Value = 1.074121072107445
$ws.Range("L21").Value = 1.060362744333984
$ws.Range("M21").Value = 1.07994033777757
$ws.Range("N21").Value = 1.069185138881287
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.06002623160189
$ws.Range("D22").Value = 1.06947813078802
$ws.Range("E22").Value = 1.055788810060623
$ws.Range("F22").Value = 1.075288486481994
$ws.Range("I22").Value = 1.043356915983291
$ws.Range("J22").Value = 1.066680590838228
$ws.Range("K22").Value = 1.07306234117604
$ws.Range("L22").Value = 1.059423036721592
$ws.Range("M22").Value = 1.078851892550785
$ws.Range("N22").Value = 1.068195399627193
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.060680547250201
$ws.Range("D23").Value = 1.070107759744592
$ws.Range("E23").Value = 1.05635614551876
$ws.Range("F23").Value = 1.075933484860117
$ws.Range("I23").Value = 1.043513570063298
$ws.Range("J23").Value = 1.067204852380672
$ws.Range("K23").Value = 1.073623905646498
$ws.Range("L23").Value = 1.059921491938012
$ws.Range("M23").Value = 1.079429205327726
$ws.Range("N23").Value = 1.06872040568116
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.063249925800822
$ws.Range("D24").Value = 1.072580728581044
$ws.Range("E24").Value = 1.058584218123172
$ws.Range("F24").Value = 1.078467033880601
$ws.Range("I24").Value = 1.044124503564225
$ws.Range("J24").Value = 1.06926166455445
$ws.Range("K24").Value = 1.075827910522364
$ws.Range("L24").Value = 1.061877315150931
$ws.Range("M24").Value = 1.081695295750829
$ws.Range("N24").Value = 1.070780138764145
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.066218454481262
$ws.Range("D25").Value = 1.075438939522638
$ws.Range("E25").Value = 1.061158897926343
$ws.Range("F25").Value = 1.081395705090347
$ws.Range("I25").Value = 1.04482160747033
$ws.Range("J25").Value = 1.071634121299006
$ws.Range("K25").Value = 1.078371892935471
$ws.Range("L25").Value = 1.064133807514384
$ws.Range("M25").Value = 1.08638783758129
$ws.Range("N25").Value = 1.073155964669403
